$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

$ws.Range("A$row").Value = "2ac21e18-91a2-45af-b5b2-dadf7b7d688c"
$ws.Range("B$row").Value = "dayli"
$ws.Range("C$row").Value = "Juices"
$ws.Range("D$row").Value = "Coca-Cola"
$ws.Range("E$row").Value = 95
$ws.Range("F$row").Value = 30

# Store Date/Time as plain text (matching the other rows), not as
# auto-converted date/time serial numbers.
$ws.Range("G$row").NumberFormat = "@"
$ws.Range("G$row").Value = "2024-09-20"
$ws.Range("H$row").NumberFormat = "@"
$ws.Range("H$row").Value = "18:22:11"
